# Apply FHIR ValueSet metadata refresh (5.0.0 -> 6.0.0) to the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Publication date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$ws.Range("B9").Value = "Alvearie Team"

# Replace the two duplicated "Contact" rows with Jurisdiction and Description
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "IBM Health Data Connect standard values for practitioner role type"

# Shift remaining property rows down by one
$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").ClearContents()

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").ClearContents()

$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# Remove the now-duplicate trailing Immutable row
$ws.Rows.Item(15).Delete()
